$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds the last-changed date serial for each data row.
# Bump every value in C2:C176 from 45204 to 45205 (one day later), matching
# the automatic-update diff which changed every row's C cell by +1.
$ws.Range("C2:C176").Value = 45205
